$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D3").Value = "fjkldfj"
$ws.Activate()
$ws.Range("D3").Select()
